$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 327.14285
$ws.Range("I19").Value = 429.33334
$ws.Range("J19").Value = 250.5
$ws.Range("K19").Value = 429.33334
$ws.Range("L19").Value = 250.5
$ws.Range("M19").Value = -254.33334
$ws.Range("N19").Value = -600.5
$ws.Range("H33").Value = 377.54544
$ws.Range("I33").Value = 405.5
$ws.Range("K33").Value = 405.5
$ws.Range("M33").Value = -176.5
$ws.Range("H40").Value = 2600
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825
$ws.Range("H70").Value = 3590.818
$ws.Range("I70").Value = 3388.889
$ws.Range("K70").Value = 10166.667
$ws.Range("M70").Value = -9896.667000000001
$ws.Range("H73").Value = 3590.818
$ws.Range("I73").Value = 3388.889
$ws.Range("K73").Value = 10166.667
$ws.Range("M73").Value = -9230.667000000001
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H92").Value = 50000800
$ws.Range("I92").Value = 76923896
$ws.Range("J92").Value = 768.1429000000001
$ws.Range("K92").Value = 76923896
$ws.Range("L92").Value = 768.1429000000001
$ws.Range("M92").Value = -76922648
$ws.Range("N92").Value = -3264.1429
$ws.Range("H116").Value = 5000
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884
$ws.Range("H132").Value = 2010.5333
$ws.Range("I132").Value = 2118.4285
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 6355.2855
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -3825.2855
$ws.Range("N132").Value = -6560
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 4756.6113
$ws.Range("I138").Value = 4187.125
$ws.Range("J138").Value = 4919.3213
$ws.Range("K138").Value = 12561.375
$ws.Range("L138").Value = 14757.9639
$ws.Range("M138").Value = -7421.375
$ws.Range("N138").Value = -25037.9639

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 600
$ws.Range("J5").Value = 600
$ws.Range("L5").Value = 600
$ws.Range("N5").Value = -824
$ws.Range("H38").Value = 4254
$ws.Range("I38").Value = 3812.25
$ws.Range("K38").Value = 3812.25
$ws.Range("M38").Value = -3345.25
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H88").Value = 2616.5
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 2674.75
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 2674.75
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -3486.75
$ws.Range("H91").Value = 2616.5
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 2674.75
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 2674.75
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -5482.75
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 600
$ws.Range("J4").Value = 600
$ws.Range("L4").Value = 600
$ws.Range("N4").Value = -830
$ws.Range("H62").Value = 130000
$ws.Range("J62").Value = 130000
$ws.Range("L62").Value = 130000
$ws.Range("N62").Value = -131372
$ws.Range("H65").Value = 130000
$ws.Range("J65").Value = 130000
$ws.Range("L65").Value = 390000
$ws.Range("N65").Value = -396864
$ws.Range("H99").Value = 2861.1428
$ws.Range("I99").Value = 2861.1428
$ws.Range("K99").Value = 2861.1428
$ws.Range("M99").Value = -1363.1428

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9164.833000000001
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 9164.833000000001
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 9164.833000000001
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -9754.833000000001
$ws.Range("H34").Value = 9164.833000000001
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9164.833000000001
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 9164.833000000001
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -9568.833000000001
$ws.Range("H35").Value = 992.4
$ws.Range("I35").Value = 992.4
$ws.Range("K35").Value = 992.4
$ws.Range("M35").Value = -698.4
$ws.Range("H58").Value = 900
$ws.Range("I58").Value = 900
$ws.Range("K58").Value = 900
$ws.Range("M58").Value = -697
$ws.Range("H134").Value = 2000
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465
$ws.Range("H136").Value = 900
$ws.Range("I136").Value = 900
$ws.Range("K136").Value = 2700
$ws.Range("M136").Value = -150
$ws.Range("H141").Value = 813610
$ws.Range("J141").Value = 813610
$ws.Range("L141").Value = 813610
$ws.Range("N141").Value = -823970

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 766.6
$ws.Range("J113").Value = 708.25
$ws.Range("L113").Value = 2124.75
$ws.Range("N113").Value = -6464.75
$ws.Range("H115").Value = 2256.5
$ws.Range("I115").Value = 1342
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 4026
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = -2851
$ws.Range("N115").Value = -17350

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1710
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H102").Value = 1914
$ws.Range("I102").Value = 1225
$ws.Range("K102").Value = 1225
$ws.Range("M102").Value = 397
$ws.Range("H104").Value = 27000
$ws.Range("J104").Value = 27000
$ws.Range("L104").Value = 27000
$ws.Range("N104").Value = -33988
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 8466.333000000001
$ws.Range("I20").Value = 8466.333000000001
$ws.Range("K20").Value = 8466.333000000001
$ws.Range("M20").Value = -8240.333000000001
$ws.Range("H22").Value = 693.75
$ws.Range("I22").Value = 387.5
$ws.Range("K22").Value = 387.5
$ws.Range("M22").Value = -92.5
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 693.75
$ws.Range("I27").Value = 387.5
$ws.Range("K27").Value = 387.5
$ws.Range("M27").Value = -280.5
$ws.Range("H46").Value = 897.5
$ws.Range("I46").Value = 897.5
$ws.Range("K46").Value = 897.5
$ws.Range("M46").Value = -709.5
$ws.Range("H136").Value = 1999.6666
$ws.Range("I136").Value = 1999.6666
$ws.Range("K136").Value = 5998.9998
$ws.Range("M136").Value = -3448.9998

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 3000
$ws.Range("I29").Value = 3000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2710
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("H62").Value = 6998.75
$ws.Range("I62").Value = 6333.3335
$ws.Range("J62").Value = 8995
$ws.Range("K62").Value = 6333.3335
$ws.Range("L62").Value = 8995
$ws.Range("M62").Value = -5709.3335
$ws.Range("N62").Value = -10243
$ws.Range("H65").Value = 6998.75
$ws.Range("I65").Value = 6333.3335
$ws.Range("J65").Value = 8995
$ws.Range("K65").Value = 31666.6675
$ws.Range("L65").Value = 44975
$ws.Range("M65").Value = -28546.6675
$ws.Range("N65").Value = -51215
$ws.Range("H107").Value = 317.25
$ws.Range("I107").Value = 293.8
$ws.Range("J107").Value = 356.33334
$ws.Range("K107").Value = 881.4000000000001
$ws.Range("L107").Value = 1069.00002
$ws.Range("M107").Value = 1038.6
$ws.Range("N107").Value = -4909.000019999999
$ws.Range("H122").Value = 1656.75
$ws.Range("I122").Value = 1575.7858
$ws.Range("K122").Value = 4727.357400000001
$ws.Range("M122").Value = -2277.357400000001
